# Actualizacion Datos Personales 4 nov
# Insert a new student record into the "Rescatables" sheet, as row 4,
# pushing the existing rows 4-10 down to rows 5-11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Insert a new blank row at position 4 (shifts rows 4..10 down to 5..11)
$ws.Rows.Item(4).Insert()

# Fill in the new student's data
$ws.Cells.Item(4, 1).Value = 20330051920081
$ws.Cells.Item(4, 2).Value = "GONZALEZ"
$ws.Cells.Item(4, 3).Value = "MENDEZ"
$ws.Cells.Item(4, 4).Value = "CRISTIAN JAHIR"
$ws.Cells.Item(4, 5).Value = "BIOLOGÍA"
$ws.Cells.Item(4, 6).Value = "3AEV"
$ws.Cells.Item(4, 7).Value = 6
